$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the descriptions to reflect payment/delivery now handled via Flag
$ws.Range("C17").Value = "Lista de Pedidos Pagos (através de Flag)"
$ws.Range("C16").Value = "Lista de Pedidos Entregues (através de Flag)"

# Move the active selection to A5
$ws.Range("A5").Select()
